$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.339.31'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.568.25'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.89'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.22'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '23.79'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.30%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0586'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0895'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.794.35'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.575.64'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.66'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.344.66'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.512'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.07'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '227.15'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.37'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.67%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.95%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.92'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.59%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.89%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.75'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.87'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.57%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.73%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0480'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.72%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.88%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.377.63'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.31%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.90%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.98%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.86%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.09%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.98%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.24%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0473'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.780'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.72%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '62.15'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.02%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -6.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.706.53'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '85.38'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.75%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.97%  '
